# The upstream commit ("Fixed #295 Add the version of M2Doc in the
# template custom properties.") touches many files in the repository,
# but the OOXML diff recorded for *this* particular template
# (notExistingLogin-template.docx) contains no semantic change at all:
# every "-"/"+" line pair is the exact same element/attribute set, just
# re-serialized with a different (alphabetical) attribute order, e.g.
#   <w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/>
#   -> <w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/>
# No text, field code, style value, page size/margin, font, language,
# or custom property actually changes anywhere in the diff (verified
# attribute-set equality for all 41 changed elements). That reordering
# is a byproduct of whatever XML (re)serializer produced the commit's
# snapshot, not an addressable Word object-model edit - Word's COM
# automation surface has no "attribute order" concept to control, and
# this document's content/formatting already matches the target state.
#
# So there is nothing to edit here: leave the document's content and
# formatting exactly as-is.
$d = $word.ActiveDocument
